# Auto-generated edit script: applies cached-value corrections
# to the Goblin_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the authoritative diff of the workbook's OOXML.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8891.857
$ws.Range("I62").Value = 2996
$ws.Range("J62").Value = 12167.333
$ws.Range("K62").Value = 2996
$ws.Range("L62").Value = 12167.333
$ws.Range("M62").Value = -2372
$ws.Range("N62").Value = -13415.333
$ws.Range("H65").Value = 8891.857
$ws.Range("I65").Value = 2996
$ws.Range("J65").Value = 12167.333
$ws.Range("K65").Value = 14980
$ws.Range("L65").Value = 60836.665
$ws.Range("M65").Value = -11860
$ws.Range("N65").Value = -67076.66500000001
$ws.Range("H92").Value = 409.83334
$ws.Range("I92").Value = 472.08334
$ws.Range("J92").Value = 160.83333
$ws.Range("K92").Value = 472.08334
$ws.Range("L92").Value = 160.83333
$ws.Range("M92").Value = 775.91666
$ws.Range("N92").Value = -2656.83333
$ws.Range("H99").Value = 1094.8182
$ws.Range("I99").Value = 704.3
$ws.Range("K99").Value = 2112.9
$ws.Range("M99").Value = -614.8999999999996
$ws.Range("H100").Value = 5587.1177
$ws.Range("I100").Value = 3998.1428
$ws.Range("J100").Value = 6699.4
$ws.Range("K100").Value = 3998.1428
$ws.Range("L100").Value = 6699.4
$ws.Range("M100").Value = -3457.1428
$ws.Range("N100").Value = -7781.4
$ws.Range("H103").Value = 1593.3334
$ws.Range("J103").Value = 1500
$ws.Range("L103").Value = 4500
$ws.Range("N103").Value = -5672
$ws.Range("H106").Value = 6375
$ws.Range("J106").Value = 4000
$ws.Range("L106").Value = 4000
$ws.Range("N106").Value = -5262
$ws.Range("H113").Value = 12799.875
$ws.Range("I113").Value = 24000
$ws.Range("K113").Value = 24000
$ws.Range("M113").Value = -20746
$ws.Range("H137").Value = 2154
$ws.Range("I137").Value = 2647.625
$ws.Range("J137").Value = 1166.75
$ws.Range("K137").Value = 7942.875
$ws.Range("L137").Value = 3500.25
$ws.Range("M137").Value = -5392.875
$ws.Range("N137").Value = -8600.25
$ws.Range("H138").Value = 4695.44
$ws.Range("I138").Value = 4138.6
$ws.Range("J138").Value = 4834.65
$ws.Range("K138").Value = 12415.8
$ws.Range("L138").Value = 14503.95
$ws.Range("M138").Value = -7275.800000000001
$ws.Range("N138").Value = -24783.95
$ws.Range("H141").Value = 4857
$ws.Range("I141").Value = 4799.8
$ws.Range("K141").Value = 14399.4
$ws.Range("M141").Value = -9219.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1058.8334
$ws.Range("I45").Value = 967.8
$ws.Range("K45").Value = 967.8
$ws.Range("M45").Value = -590.8
$ws.Range("H97").Value = 788.25
$ws.Range("I97").Value = 811.5833
$ws.Range("J97").Value = 718.25
$ws.Range("K97").Value = 811.5833
$ws.Range("L97").Value = 718.25
$ws.Range("M97").Value = -315.5833
$ws.Range("N97").Value = -1710.25
$ws.Range("H102").Value = 3141.3225
$ws.Range("I102").Value = 1822.3462
$ws.Range("K102").Value = 1822.3462
$ws.Range("M102").Value = -200.3462
$ws.Range("H109").Value = 120125.664
$ws.Range("J109").Value = 120125.664
$ws.Range("L109").Value = 120125.664
$ws.Range("N109").Value = -122899.664
$ws.Range("H110").Value = 1200
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H131").Value = 140000
$ws.Range("J131").Value = 140000
$ws.Range("L131").Value = 140000
$ws.Range("N131").Value = -150080
$ws.Range("H132").Value = 2353.547
$ws.Range("I132").Value = 2351.898
$ws.Range("J132").Value = 2373.75
$ws.Range("K132").Value = 7055.694
$ws.Range("L132").Value = 7121.25
$ws.Range("M132").Value = -4525.694
$ws.Range("N132").Value = -12181.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 348.25
$ws.Range("J64").Value = 348.25
$ws.Range("L64").Value = 348.25
$ws.Range("N64").Value = -798.25
$ws.Range("H67").Value = 348.25
$ws.Range("J67").Value = 348.25
$ws.Range("L67").Value = 348.25
$ws.Range("N67").Value = -1908.25
$ws.Range("H94").Value = 2059.3215
$ws.Range("I94").Value = 1835.875
$ws.Range("K94").Value = 1835.875
$ws.Range("M94").Value = -1384.875
$ws.Range("H99").Value = 2888.2856
$ws.Range("I99").Value = 1524.5
$ws.Range("J99").Value = 4706.6665
$ws.Range("K99").Value = 1524.5
$ws.Range("L99").Value = 4706.6665
$ws.Range("M99").Value = -26.5
$ws.Range("N99").Value = -7702.6665
$ws.Range("H120").Value = 74000
$ws.Range("J120").Value = 74000
$ws.Range("L120").Value = 74000
$ws.Range("N120").Value = -83676
$ws.Range("H134").Value = 3658.182
$ws.Range("I134").Value = 3842.625
$ws.Range("J134").Value = 3166.3333
$ws.Range("K134").Value = 11527.875
$ws.Range("L134").Value = 9498.999899999999
$ws.Range("M134").Value = -8992.875
$ws.Range("N134").Value = -14568.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 4566.3335
$ws.Range("I37").Value = 3999
$ws.Range("J37").Value = 4850
$ws.Range("K37").Value = 3999
$ws.Range("L37").Value = 4850
$ws.Range("M37").Value = -3892
$ws.Range("N37").Value = -5064
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H99").Value = 2960
$ws.Range("I99").Value = 2552
$ws.Range("K99").Value = 2552
$ws.Range("M99").Value = -1054
$ws.Range("H107").Value = 1720
$ws.Range("J107").Value = 1399
$ws.Range("L107").Value = 1399
$ws.Range("N107").Value = -5239
$ws.Range("H126").Value = 2960
$ws.Range("I126").Value = 2552
$ws.Range("K126").Value = 7656
$ws.Range("M126").Value = -5186
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2932
$ws.Range("I5").Value = 1062
$ws.Range("K5").Value = 3186
$ws.Range("M5").Value = -3074
$ws.Range("H12").Value = 105.666664
$ws.Range("I12").Value = 17.75
$ws.Range("K12").Value = 53.25
$ws.Range("M12").Value = 119.75
$ws.Range("H51").Value = 1205.4445
$ws.Range("I51").Value = 274.5
$ws.Range("J51").Value = 1471.4286
$ws.Range("K51").Value = 823.5
$ws.Range("L51").Value = 4414.2858
$ws.Range("M51").Value = -363.5
$ws.Range("N51").Value = -5334.2858
$ws.Range("H62").Value = 2415.4
$ws.Range("I62").Value = 756
$ws.Range("J62").Value = 3521.6667
$ws.Range("K62").Value = 2268
$ws.Range("L62").Value = 10565.0001
$ws.Range("M62").Value = -1582
$ws.Range("N62").Value = -11937.0001
$ws.Range("H65").Value = 2415.4
$ws.Range("I65").Value = 756
$ws.Range("J65").Value = 3521.6667
$ws.Range("K65").Value = 6804
$ws.Range("L65").Value = 31695.0003
$ws.Range("M65").Value = -3372
$ws.Range("N65").Value = -38559.0003
$ws.Range("H107").Value = 3126.7273
$ws.Range("J107").Value = 2413.8333
$ws.Range("L107").Value = 7241.499899999999
$ws.Range("N107").Value = -11081.4999
$ws.Range("H130").Value = 5384.2
$ws.Range("I130").Value = 4231.5
$ws.Range("K130").Value = 12694.5
$ws.Range("M130").Value = -7674.5
$ws.Range("H135").Value = 2932
$ws.Range("I135").Value = 1062
$ws.Range("K135").Value = 9558
$ws.Range("M135").Value = -7023
$ws.Range("H140").Value = 4850.6665
$ws.Range("I140").Value = 3821
$ws.Range("K140").Value = 11463
$ws.Range("M140").Value = -6283

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4544.7856
$ws.Range("I97").Value = 1126.6818
$ws.Range("K97").Value = 1126.6818
$ws.Range("M97").Value = -630.6818000000001
$ws.Range("H113").Value = 7417.591
$ws.Range("I113").Value = 1324.2
$ws.Range("K113").Value = 1324.2
$ws.Range("M113").Value = 845.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4865.1304
$ws.Range("I93").Value = 4133.3335
$ws.Range("K93").Value = 4133.3335
$ws.Range("M93").Value = -2885.3335
$ws.Range("H136").Value = 5750.25
$ws.Range("I136").Value = 5143.2856
$ws.Range("K136").Value = 15429.8568
$ws.Range("M136").Value = -12879.8568
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5098.6313
$ws.Range("I96").Value = 6180.091
$ws.Range("J96").Value = 3611.625
$ws.Range("K96").Value = 6180.091
$ws.Range("L96").Value = 3611.625
$ws.Range("M96").Value = -4807.091
$ws.Range("N96").Value = -6357.625
$ws.Range("H100").Value = 3000.2856
$ws.Range("I100").Value = 2599.8
$ws.Range("J100").Value = 4001.5
$ws.Range("K100").Value = 5199.6
$ws.Range("L100").Value = 8003
$ws.Range("M100").Value = -4658.6
$ws.Range("N100").Value = -9085
$ws.Range("H137").Value = 57856.57
$ws.Range("J137").Value = 57856.57
$ws.Range("L137").Value = 57856.57
$ws.Range("N137").Value = -68056.57000000001

